$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "heuristica" evaluation row (row 4) to the "Heurística" sheet
#    and the new shared string "A" used by AU4.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Heurística")

$ws3.Range("E4").Value  = 5
$ws3.Range("F4").Value  = 4
$ws3.Range("G4").Value  = 4.5
$ws3.Range("H4").Value  = 5
$ws3.Range("I4").Value  = 3.5
$ws3.Range("J4").Value  = 5
$ws3.Range("K4").Value  = 3
$ws3.Range("L4").Value  = 1
$ws3.Range("M4").Value  = 4
$ws3.Range("N4").Value  = 5
$ws3.Range("O4").Value  = 5
$ws3.Range("P4").Value  = 5
$ws3.Range("Q4").Value  = 3.5
$ws3.Range("R4").Value  = 3
$ws3.Range("S4").Value  = 4
$ws3.Range("T4").Value  = 4
$ws3.Range("U4").Value  = 4
$ws3.Range("V4").Value  = 4
$ws3.Range("W4").Value  = 4.5
$ws3.Range("X4").Value  = 5
$ws3.Range("Y4").Value  = 4.5
$ws3.Range("Z4").Value  = 5
$ws3.Range("AA4").Value = 3.5
$ws3.Range("AB4").Value = 5
$ws3.Range("AC4").Value = 5
$ws3.Range("AD4").Value = 5
$ws3.Range("AE4").Value = 5
$ws3.Range("AF4").Value = 5
$ws3.Range("AG4").Value = 4.5
$ws3.Range("AH4").Value = 5
$ws3.Range("AI4").Value = 5
$ws3.Range("AJ4").Value = 4.5
$ws3.Range("AK4").Value = 5
$ws3.Range("AL4").Value = 5
$ws3.Range("AM4").Value = 5
$ws3.Range("AN4").Value = 5
$ws3.Range("AO4").Value = 5
$ws3.Range("AP4").Value = 5
$ws3.Range("AQ4").Value = 5
$ws3.Range("AR4").Value = 5
$ws3.Range("AS4").Value = 3
$ws3.Range("AT4").Value = 3.5
$ws3.Range("AU4").Value = "A"
$ws3.Range("AV4").Value = 0
$ws3.Range("AW4").Value = 5
$ws3.Range("AX4").Value = 5

# ---------------------------------------------------------------------------
# 2. Update the saved view/selection state on the various "Impacto" sheets.
#    Selecting a range on a sheet also activates that sheet (as in real
#    Excel), so we do these before finally activating "Heurística" so the
#    workbook ends up with that sheet as the active tab, matching the diff.
# ---------------------------------------------------------------------------
$wsS1 = $wb.Worksheets.Item("Impacto S1")
[void]$wsS1.Range("E42:E49").Select()

$wsS7 = $wb.Worksheets.Item("Impacto S7")
[void]$wsS7.Range("D43:D50").Select()

$wsS8 = $wb.Worksheets.Item("Impacto S8")
[void]$wsS8.Range("D43:D50").Select()

$wsS9 = $wb.Worksheets.Item("Impacto S9")
[void]$wsS9.Range("D43:D50").Select()

$wsS10 = $wb.Worksheets.Item("Impacto S10")
[void]$wsS10.Range("D43:E50").Select()

# ---------------------------------------------------------------------------
# 3. Finally select/activate the new active cell on "Heurística" so that it
#    becomes the active tab (activeTab goes from 4 -> 2) and its selection
#    becomes AY4.
# ---------------------------------------------------------------------------
[void]$ws3.Range("AY4").Select()

Write-Output "edit applied"
